# Glenn Phillips.xlsx — add match-level rows + a "matchNo" column,
# and rename the sheet/tab to the player's name.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename the sheet/tab to match the player name -------------------------
$ws.Name = "Glenn Phillips"

# --- Insert a new first column "matchNo" (shifts teamName..result right) ---
$ws.Columns.Item(1).Insert()

# --- Insert two more match rows above the existing (only) data row ---------
# The single pre-existing data row (b Shivam Mavi / KKR match) becomes the
# third data row (row 4) once the two new ones are inserted above it.
$ws.Rows.Item(2).Insert()
$ws.Rows.Item(2).Insert()

# Every value in this sheet is stored as text (even the numeric-looking
# ones), matching the original file's convention. Force text storage with a
# leading quote-prefix for cells that would otherwise be auto-parsed as a
# number by Excel.

# --- Header row --------------------------------------------------------
$ws.Range("A1").Value = "matchNo"
$ws.Range("B1").Value = "teamName"
$ws.Range("C1").Value = "batterName"
$ws.Range("D1").Value = "states"
$ws.Range("E1").Value = "runs"
$ws.Range("F1").Value = "balls"
$ws.Range("G1").Value = "fours"
$ws.Range("H1").Value = "sixes"
$ws.Range("I1").Value = "sr"
$ws.Range("J1").Value = "opponentTeamName"
$ws.Range("K1").Value = "venue"
$ws.Range("L1").Value = "date"
$ws.Range("M1").Value = "result"

# --- Row 2: 47th match vs Chennai Super Kings ---------------------------
$ws.Range("A2").Value = "47th"
$ws.Range("B2").Value = "Rajasthan Royals"
$ws.Range("C2").Value = "Glenn Phillips"
$ws.Range("D2").Value = ""
$ws.Range("E2").Value = "'14"
$ws.Range("F2").Value = "'8"
$ws.Range("G2").Value = "'1"
$ws.Range("H2").Value = "'1"
$ws.Range("I2").Value = "'175.00"
$ws.Range("J2").Value = "Chennai Super Kings"
$ws.Range("K2").Value = "Abu Dhabi"
$ws.Range("L2").Value = "October 02"
$ws.Range("M2").Value = "Royals won by 7 wickets (with 15 balls remaining)"

# --- Row 3: 51st match vs Mumbai Indians --------------------------------
$ws.Range("A3").Value = "51st"
$ws.Range("B3").Value = "Rajasthan Royals"
$ws.Range("C3").Value = "Glenn Phillips"
$ws.Range("D3").Value = "b Coulter-Nile"
$ws.Range("E3").Value = "'4"
$ws.Range("F3").Value = "'13"
$ws.Range("G3").Value = "'0"
$ws.Range("H3").Value = "'0"
$ws.Range("I3").Value = "'30.76"
$ws.Range("J3").Value = "Mumbai Indians"
$ws.Range("K3").Value = "Sharjah"
$ws.Range("L3").Value = "October 05"
$ws.Range("M3").Value = "Mumbai won by 8 wickets (with 70 balls remaining)"

# --- Row 4: 54th match vs Kolkata Knight Riders (pre-existing row data) -
$ws.Range("A4").Value = "54th"
$ws.Range("B4").Value = "Rajasthan Royals"
$ws.Range("C4").Value = "Glenn Phillips"
$ws.Range("D4").Value = "b Shivam Mavi"
$ws.Range("E4").Value = "'8"
$ws.Range("F4").Value = "'12"
$ws.Range("G4").Value = "'0"
$ws.Range("H4").Value = "'1"
$ws.Range("I4").Value = "'66.66"
$ws.Range("J4").Value = "Kolkata Knight Riders"
$ws.Range("K4").Value = "Sharjah"
$ws.Range("L4").Value = "October 07"
$ws.Range("M4").Value = "KKR won by 86 runs"
